$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Version: 5.0.0 -> 6.0.0
$ws.Range("B3").Value = "6.0.0"

# Date: 2021-12-16T17:36:56+00:00 -> 2022-01-21T20:46:54+00:00
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value (was blank) -> Alvearie Team
$ws.Range("B9").Value = "Alvearie Team"

# Row 10: Contact / No display for ContactDetail -> Jurisdiction / United States of America
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# Remove the now-duplicate old "Contact" row entirely, shifting everything below up by one
$ws.Rows.Item(11).Delete()
